$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$ref, [string]$val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue "D2" '43.501.65'
Set-TextValue "E2" '  +4.53%  '
Set-TextValue "D3" '2.273.92'
Set-TextValue "E3" '  +2.38%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.12%  '
Set-TextValue "D5" '231.67'
Set-TextValue "E5" '  +0.58%  '
Set-TextValue "D6" '0.629'
Set-TextValue "E6" '  +0.91%  '
Set-TextValue "D7" '61.43'
Set-TextValue "E7" '  +0.73%  '
Set-TextValue "E8" '  -0.01%  '
Set-TextValue "D9" '0.413'
Set-TextValue "E9" '  +2.92%  '
Set-TextValue "D10" '0.0918'
Set-TextValue "E10" '  +3.11%  '
Set-TextValue "E11" '  +0.42%  '
Set-TextValue "D12" '2.603.42'
Set-TextValue "E12" '  +2.10%  '
Set-TextValue "D13" '15.77'
Set-TextValue "E13" '  +0.61%  '
Set-TextValue "D14" '22.65'
Set-TextValue "E14" '  +4.27%  '
Set-TextValue "D15" '5.71'
Set-TextValue "E15" '  +2.91%  '
Set-TextValue "D16" '0.811'
Set-TextValue "E16" '  +1.73%  '
Set-TextValue "D17" '2.267.72'
Set-TextValue "E17" '  +2.11%  '
Set-TextValue "D18" '43.327.82'
Set-TextValue "E18" '  +4.32%  '
Set-TextValue "D19" '0.0₃0933'
Set-TextValue "E19" '  +4.29%  '
Set-TextValue "D20" '73.12'
Set-TextValue "E20" '  +0.55%  '
Set-TextValue "D21" '6.21'
Set-TextValue "E21" '  +2.74%  '
Set-TextValue "D22" '250.22'
Set-TextValue "E22" '  -0.27%  '
Set-TextValue "D24" '2.56'
Set-TextValue "E24" '  +7.54%  '
Set-TextValue "E25" '  +4.22%  '
Set-TextValue "D26" '9.80'
Set-TextValue "E26" '  +2.60%  '
Set-TextValue "D27" '169.73'
Set-TextValue "E27" '  +1.20%  '
Set-TextValue "E28" '  +1.35%  '
Set-TextValue "B29" 'ImmutableX'
Set-TextValue "C29" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D29" '1.50'
Set-TextValue "E29" '  +6.38%  '
Set-TextValue "B30" 'EthereumClassic'
Set-TextValue "C30" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D30" '20.62'
Set-TextValue "E30" '  +3.44%  '
Set-TextValue "D31" '2.66'
Set-TextValue "E31" '  +0.81%  '
Set-TextValue "E32" '  -0.46%  '
Set-TextValue "D33" '5.03'
Set-TextValue "E33" '  +2.14%  '
Set-TextValue "D34" '4.74'
Set-TextValue "E34" '  +2.97%  '
Set-TextValue "D35" '0.0656'
Set-TextValue "E35" '  +5.52%  '
Set-TextValue "B36" 'LidoDAOToken'
Set-TextValue "C36" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D36" '2.40'
Set-TextValue "E36" '  +2.02%  '
Set-TextValue "B37" 'THORChain'
Set-TextValue "C37" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D37" '6.46'
Set-TextValue "E37" '  -1.73%  '
Set-TextValue "D38" '3.61'
Set-TextValue "E38" '  -1.82%  '
Set-TextValue "D39" '0.0251'
Set-TextValue "E39" '  +5.33%  '
Set-TextValue "E40" '  -0.06%  '
Set-TextValue "D41" '8.69'
Set-TextValue "E41" '  +0.89%  '
Set-TextValue "D42" '0.000221'
Set-TextValue "E42" '  -9.58%  '
Set-TextValue "D43" '0.0968'
Set-TextValue "E43" '  -0.84%  '
Set-TextValue "D44" '1.21'
Set-TextValue "E44" '  +0.35%  '
Set-TextValue "B45" 'Aave'
Set-TextValue "C45" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D45" '97.33'
Set-TextValue "E45" '  -1.31%  '
Set-TextValue "B46" 'FTXToken'
Set-TextValue "C46" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D46" '4.39'
Set-TextValue "E46" '  -9.19%  '
Set-TextValue "D47" '1.465.39'
Set-TextValue "E47" '  -0.03%  '
Set-TextValue "E48" '  +1.44%  '
Set-TextValue "B49" 'ARBITRUM'
Set-TextValue "C49" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D49" '1.08'
Set-TextValue "E49" '  +1.03%  '
Set-TextValue "B50" 'HuobiToken'
Set-TextValue "C50" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D50" '2.76'
Set-TextValue "E50" '  -1.68%  '
Set-TextValue "B51" 'NEARProtocol'
Set-TextValue "C51" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D51" '2.25'
Set-TextValue "E51" '  +7.16%  '
